# Weekly update: new price-report rows arrived for "Zapallo italiano" /
# Agrícola del Norte S.A. de Arica. The source feed always keeps the most
# recent report at the top of the data block (row 146 onward), so this
# week's two new quality rows (Primera / Segunda) are inserted right above
# the existing data, pushing every older row down by two rows and growing
# the used range from A1:R228 to A1:R230.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 146-147),
# shifting the existing rows 146..228 down to 148..230.
$ws.Rows.Item(146).Resize(2).Insert()

# New row 146 - "Primera" quality, week of 2021-11-29 (serial 44529)
$newRow146 = @(
    1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44529, 15,
    100112032, "Zapallo italiano", "Huracán", "Primera", 120, 4500, 5000,
    4750, "`$/caja 70 unidades", "Región de Arica y Parinacota", 68, 70,
    "Hortaliza"
)
for ($c = 1; $c -le $newRow146.Length; $c++) {
    $ws.Cells.Item(146, $c).Value = $newRow146[$c - 1]
}

# New row 147 - "Segunda" quality, same week
$newRow147 = @(
    1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44529, 15,
    100112032, "Zapallo italiano", "Huracán", "Segunda", 120, 3500, 4000,
    3750, "`$/caja 100 unidades", "Región de Arica y Parinacota", 38, 100,
    "Hortaliza"
)
for ($c = 1; $c -le $newRow147.Length; $c++) {
    $ws.Cells.Item(147, $c).Value = $newRow147[$c - 1]
}
